$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (EMMERICH, HAL)
$ws.Range("C2").Value = 1234
$ws.Range("D2").Value = 123
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 3702
$ws.Range("H2").Value = 615
$ws.Range("I2").Value = 4317

# Update row 3 (OCELOT, REVOLVER)
$ws.Range("C3").Value = 2345
$ws.Range("D3").Value = 234
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4690
$ws.Range("H3").Value = 234
$ws.Range("I3").Value = 4924

# Update row 4 (SILVERBURGH, MERYL)
$ws.Range("C4").Value = 1234
$ws.Range("D4").Value = 123
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 3702
$ws.Range("H4").Value = 492
$ws.Range("I4").Value = 4194

# Delete rows 5-8 (SNAKE SOLID, SNAKE LIQUID, SNAKE SOLIDUS, SNAKE NAKED)
$ws.Range("A5:I8").EntireRow.Delete()
